$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Josh Green" (SG,SF / Charlotte Hornets) which is no longer present
# in the updated roster. Find it dynamically in case the row order differs.
$targetRow = 0
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($name -eq "Josh Green") {
        $targetRow = $r
        break
    }
}
if ($targetRow -gt 0) {
    $ws.Rows.Item($targetRow).Delete() | Out-Null
}

# Re-upload the roster data in its new order.
$data = @(
    @("Jose Alvarado", "PG", "New Orleans Pelicans"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Amen Thompson", "SG,SF,PF", "Houston Rockets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
